# Applies per-row updates to the cryptocurrency price/volume table.
# Price cells in column D are stored as text in the workbook, so we
# force a text number format before writing numeric-looking strings to
# avoid Excel silently converting them to floating point numbers.
# Rows 42 and 43 swap identity (Aptos <-> TrustWalletToken) with new values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.102.22"
$ws.Range("E2").Value = "  +5.60%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.918.52"
$ws.Range("E3").Value = "  +2.53%  "

$ws.Range("E4").Value = "  -0.54%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "329.93"
$ws.Range("E5").Value = "  +4.61%  "

$ws.Range("E6").Value = "  -0.53%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5221"
$ws.Range("E7").Value = "  +2.45%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4090"
$ws.Range("E8").Value = "  +4.78%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08506"
$ws.Range("E9").Value = "  +1.90%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "43.03"
$ws.Range("E10").Value = "  +3.03%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.127"
$ws.Range("E11").Value = "  +1.90%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.49"
$ws.Range("E12").Value = "  +10.35%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.422"
$ws.Range("E13").Value = "  +3.31%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.929.10"
$ws.Range("E14").Value = "  +3.08%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.426"
$ws.Range("E15").Value = "  +1.93%  "

$ws.Range("E16").Value = "  -0.62%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "95.62"
$ws.Range("E17").Value = "  +4.93%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001113"
$ws.Range("E18").Value = "  +0.97%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06723"
$ws.Range("E19").Value = "  -0.14%  "

$ws.Range("E20").Value = "  +3.15%  "

$ws.Range("E21").Value = "  -0.46%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.008"
$ws.Range("E22").Value = "  +1.45%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.103.52"
$ws.Range("E23").Value = "  +5.55%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.33"
$ws.Range("E24").Value = "  +1.88%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.223"
$ws.Range("E25").Value = "  +0.82%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.155.08"
$ws.Range("E26").Value = "  +3.46%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "160.56"
$ws.Range("E27").Value = "  -0.02%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "21.08"
$ws.Range("E28").Value = "  +2.17%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.453"
$ws.Range("E29").Value = "  +1.62%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "129.29"
$ws.Range("E30").Value = "  +1.61%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.078"
$ws.Range("E31").Value = "  +3.63%  "

$ws.Range("E32").Value = "  +1.34%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.091"
$ws.Range("E33").Value = "  +6.06%  "

$ws.Range("E34").Value = "  +0.50%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02495"
$ws.Range("E35").Value = "  +1.58%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06621"
$ws.Range("E36").Value = "  +0.65%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2210"
$ws.Range("E37").Value = "  +2.10%  "

$ws.Range("E38").Value = "  +4.29%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.192"
$ws.Range("E39").Value = "  +3.28%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.908"
$ws.Range("E40").Value = "  +0.11%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6528"
$ws.Range("E41").Value = "  +2.46%  "

$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.65"
$ws.Range("E42").Value = "  +4.99%  "

$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.247"
$ws.Range("E43").Value = "  +0.57%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6168"
$ws.Range("E44").Value = "  +2.71%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.28"
$ws.Range("E45").Value = "  +2.22%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.770"
$ws.Range("E46").Value = "  +2.37%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.082"
$ws.Range("E47").Value = "  +3.79%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.248"
$ws.Range("E48").Value = "  +2.74%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "124.62"
$ws.Range("E49").Value = "  +2.08%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.162"
$ws.Range("E50").Value = "  +10.26%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "79.85"
$ws.Range("E51").Value = "  +4.50%  "

